$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed values
$ws.Range("B2").Value = 225
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 90

# Remove row 4 entirely (was A4=1, B4=90 before the shift)
$ws.Range("A4:B4").Delete()
